# Update the "想去人数" (want-to-go count) figures in the "展览" and
# "全部类型" worksheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 16898
$ws1.Range("F6").Value  = 1649
$ws1.Range("F8").Value  = 6
$ws1.Range("F9").Value  = 392
$ws1.Range("F10").Value = 225
$ws1.Range("F11").Value = 129
$ws1.Range("F12").Value = 11715
$ws1.Range("F14").Value = 1384
$ws1.Range("F15").Value = 4653
$ws1.Range("F20").Value = 898

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 16898
$ws4.Range("F7").Value  = 1649
$ws4.Range("F9").Value  = 6
$ws4.Range("F10").Value = 392
$ws4.Range("F11").Value = 225
$ws4.Range("F12").Value = 129
$ws4.Range("F15").Value = 11715
$ws4.Range("F17").Value = 1384
$ws4.Range("F18").Value = 4653
$ws4.Range("F23").Value = 898
